$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Adam9", "Itga3", "ECs", 3, 1, 11.20127533333333, 33.603826, 0.1186573945858706, 0.1186573945858706, 3, 1, 12.77892533333333, 38.336776, 0.7206984557633653, 0.7206984557633654, 143.1402611227751, 1288.262350104976, 0.08551620104294126, 0.08551620104294128),
    @("ECs", "Adam9", "Itga3", "FAPs", 3, 1, 11.20127533333333, 33.603826, 0.1186573945858706, 0.1186573945858706, 3, 1, 0.9898276666666668, 2.969483, 0.05582372947885773, 0.05582372947885774, 11.08733222688422, 99.785990041958, 0.00662389829602772, 0.006623898296027721),
    @("ECs", "Adam9", "Itga3", "M2", 3, 1, 11.20127533333333, 33.603826, 0.1186573945858706, 0.1186573945858706, 1, 0.3333333333333333, 0.02976266666666667, 0.08928800000000001, 0.001678537697541373, 0.001678537697541373, 0.3333798239875556, 3.000418415888, 0.0001991709099044255, 0.0001991709099044255),
    @("ECs", "Adam9", "Itga3", "sCs", 3, 1, 11.20127533333333, 33.603826, 0.1186573945858706, 0.1186573945858706, 3, 1, 3.932791, 11.798373, 0.2217992770602354, 0.2217992770602354, 44.05227481945533, 396.470473375098, 0.0263181243369972, 0.0263181243369972),
    @("FAPs", "Adam9", "Itga3", "ECs", 3, 1, 47.94465366666667, 143.833961, 0.5078874966566524, 0.5078874966566524, 3, 1, 12.77892533333333, 38.336776, 0.7206984557633653, 0.7206984557633654, 612.6811493388597, 5514.130344049737, 0.3660337345419708, 0.3660337345419709),
    @("FAPs", "Adam9", "Itga3", "FAPs", 3, 1, 47.94465366666667, 143.833961, 0.5078874966566524, 0.5078874966566524, 3, 1, 0.9898276666666668, 2.969483, 0.05582372947885773, 0.05582372947885774, 47.45694466801812, 427.1125020121631, 0.02835217421905523, 0.02835217421905523),
    @("FAPs", "Adam9", "Itga3", "M2", 3, 1, 47.94465366666667, 143.833961, 0.5078874966566524, 0.5078874966566524, 1, 0.3333333333333333, 0.02976266666666667, 0.08928800000000001, 0.001678537697541373, 0.001678537697541373, 1.426960745529778, 12.842646709768, 0.0008525083092481092, 0.0008525083092481093),
    @("FAPs", "Adam9", "Itga3", "sCs", 3, 1, 47.94465366666667, 143.833961, 0.5078874966566524, 0.5078874966566524, 3, 1, 3.932791, 11.798373, 0.2217992770602354, 0.2217992770602354, 188.5563024383837, 1697.006721945453, 0.1126490795863783, 0.1126490795863783),
    @("M2", "Adam9", "Itga3", "ECs", 3, 1, 25.59984766666667, 76.799543, 0.2711844085184091, 0.2711844085184091, 3, 1, 12.77892533333333, 38.336776, 0.7206984557633653, 0.7206984557633654, 327.1385418770409, 2944.246876893368, 0.1954421844463191, 0.1954421844463191),
    @("M2", "Adam9", "Itga3", "FAPs", 3, 1, 25.59984766666667, 76.799543, 0.2711844085184091, 0.2711844085184091, 3, 1, 0.9898276666666668, 2.969483, 0.05582372947885773, 0.05582372947885774, 25.33943748291878, 228.054937346269, 0.01513852506001571, 0.01513852506001572),
    @("M2", "Adam9", "Itga3", "M2", 3, 1, 25.59984766666667, 76.799543, 0.2711844085184091, 0.2711844085184091, 1, 0.3333333333333333, 0.02976266666666667, 0.08928800000000001, 0.001678537697541373, 0.001678537697541373, 0.7619197328204445, 6.857277595384001, 0.0004551932526836096, 0.0004551932526836097),
    @("M2", "Adam9", "Itga3", "sCs", 3, 1, 25.59984766666667, 76.799543, 0.2711844085184091, 0.2711844085184091, 3, 1, 3.932791, 11.798373, 0.2217992770602354, 0.2217992770602354, 100.6788505048377, 906.109654543539, 0.06014850575939069, 0.0601485057593907),
    @("sCs", "Adam9", "Itga3", "ECs", 3, 1, 9.654369000000001, 28.963107, 0.1022707002390678, 0.1022707002390678, 3, 1, 12.77892533333333, 38.336776, 0.7206984557633653, 0.7206984557633654, 123.372460591448, 1110.352145323032, 0.07370633573213418, 0.0737063357321342),
    @("sCs", "Adam9", "Itga3", "FAPs", 3, 1, 9.654369000000001, 28.963107, 0.1022707002390678, 0.1022707002390678, 3, 1, 0.9898276666666668, 2.969483, 0.05582372947885773, 0.05582372947885774, 9.556161540409002, 86.00545386368101, 0.005709131903759071, 0.005709131903759071),
    @("sCs", "Adam9", "Itga3", "M2", 3, 1, 9.654369000000001, 28.963107, 0.1022707002390678, 0.1022707002390678, 1, 0.3333333333333333, 0.02976266666666667, 0.08928800000000001, 0.001678537697541373, 0.001678537697541373, 0.2873397664240001, 2.586057897816, 0.0001716652257052288, 0.0001716652257052288),
    @("sCs", "Adam9", "Itga3", "sCs", 3, 1, 9.654369000000001, 28.963107, 0.1022707002390678, 0.1022707002390678, 3, 1, 3.932791, 11.798373, 0.2217992770602354, 0.2217992770602354, 37.968615513879, 341.717539624911, 0.02268356737746928, 0.02268356737746928)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $row[$c]
    }
}

Write-Output "Updated $($data.Length) rows"
